$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.905473712059812
$ws.Range("D2").Value = 3.905448924006797
$ws.Range("E2").Value = 40.45309662546128
$ws.Range("F2").Value = 14.52525452576912
$ws.Range("G2").Value = 12.44752704532048
$ws.Range("H2").Value = 9.904720290307917
$ws.Range("I2").Value = 12.75548454980884
$ws.Range("M2").Value = 57.84003720018883
$ws.Range("O2").Value = 12.85841554470193
$ws.Range("C3").Value = 3.764852942667885
$ws.Range("D3").Value = 3.789696546714558
$ws.Range("E3").Value = 37.73119929783729
$ws.Range("F3").Value = 14.68125214438379
$ws.Range("G3").Value = 12.59179241730543
$ws.Range("H3").Value = 10.03823297957386
$ws.Range("I3").Value = 12.92677391096765
$ws.Range("M3").Value = 54.2573156657387
$ws.Range("O3").Value = 13.08605440471323
$ws.Range("C4").Value = 3.675661770934629
$ws.Range("D4").Value = 3.716366721210528
$ws.Range("E4").Value = 35.95960951284731
$ws.Range("F4").Value = 14.79136246438105
$ws.Range("G4").Value = 12.70999339414944
$ws.Range("H4").Value = 10.12447106022344
$ws.Range("I4").Value = 13.04144779435948
$ws.Range("M4").Value = 51.92864706475579
$ws.Range("O4").Value = 13.23583110130845
$ws.Range("C5").Value = 3.638645574814985
$ws.Range("D5").Value = 3.685946401379057
$ws.Range("E5").Value = 35.21230098746195
$ws.Range("F5").Value = 14.83971593511907
$ws.Range("G5").Value = 12.7652083651362
$ws.Range("H5").Value = 10.16067609924295
$ws.Range("I5").Value = 13.09050458425047
$ws.Range("M5").Value = 50.94729169483841
$ws.Range("O5").Value = 13.29931705816873
$ws.Range("C6").Value = 3.632460054008028
$ws.Range("D6").Value = 3.680863576838231
$ws.Range("E6").Value = 35.08667465392342
$ws.Range("F6").Value = 14.84795185878711
$ws.Range("G6").Value = 12.77479052390807
$ws.Range("H6").Value = 10.16675182533773
$ws.Range("I6").Value = 13.09878913471534
$ws.Range("M6").Value = 50.78238217002004
$ws.Range("O6").Value = 13.31000504363131
$ws.Range("C7").Value = 3.675165205415698
$ws.Range("D7").Value = 3.715958598189972
$ws.Range("E7").Value = 35.94963395749275
$ws.Range("F7").Value = 14.7920006324524
$ws.Range("G7").Value = 12.71071004502947
$ws.Range("H7").Value = 10.12495504428802
$ws.Range("I7").Value = 13.04210005376081
$ws.Range("M7").Value = 51.91554333318729
$ws.Range("O7").Value = 13.23667745527838
$ws.Range("C8").Value = 3.857605329587284
$ws.Range("D8").Value = 3.866023003388952
$ws.Range("E8").Value = 39.53532264187566
$ws.Range("F8").Value = 14.57599619573893
$ws.Range("G8").Value = 12.4908836634861
$ws.Range("H8").Value = 9.949865519351617
$ws.Range("I8").Value = 12.81253439787743
$ws.Range("M8").Value = 56.63144614659222
$ws.Range("O8").Value = 12.93478823267061
$ws.Range("C9").Value = 4.190980020149619
$ws.Range("D9").Value = 4.141227629073909
$ws.Range("E9").Value = 45.77755605998806
$ws.Range("F9").Value = 14.27152322778262
$ws.Range("G9").Value = 12.31191213287756
$ws.Range("H9").Value = 9.64075606805998
$ws.Range("I9").Value = 12.44069930196613
$ws.Range("M9").Value = 64.85909885559037
$ws.Range("O9").Value = 12.4253217240489
$ws.Range("C10").Value = 4.41899598454978
$ws.Range("D10").Value = 4.330444601490361
$ws.Range("E10").Value = 49.89247115506764
$ws.Range("F10").Value = 14.12791004217895
$ws.Range("G10").Value = 12.35524804906226
$ws.Range("H10").Value = 9.435171739880143
$ws.Range("I10").Value = 12.21952533631949
$ws.Range("M10").Value = 70.28760507029757
$ws.Range("O10").Value = 12.10603144985317
$ws.Range("C11").Value = 4.51870931586717
$ws.Range("D11").Value = 4.413475864618293
$ws.Range("E11").Value = 51.66438077649066
$ws.Range("F11").Value = 14.08156126745253
$ws.Range("G11").Value = 12.41657632554686
$ws.Range("H11").Value = 9.346477326387436
$ws.Range("I11").Value = 12.13118059189596
$ws.Range("M11").Value = 72.62507423863961
$ws.Range("O11").Value = 11.97386099466021
$ws.Range("C12").Value = 4.555867656098341
$ws.Range("D12").Value = 4.444463868016824
$ws.Range("E12").Value = 52.32116809572973
$ws.Range("F12").Value = 14.06686877974704
$ws.Range("G12").Value = 12.4460078160457
$ws.Range("H12").Value = 9.313600211445939
$ws.Range("I12").Value = 12.0995762663881
$ws.Range("M12").Value = 73.49139680176123
$ws.Range("O12").Value = 11.92579554081342
$ws.Range("C13").Value = 4.547892014516262
$ws.Range("D13").Value = 4.437810480342423
$ws.Range("E13").Value = 52.18034563227636
$ws.Range("F13").Value = 14.06990389740248
$ws.Range("G13").Value = 12.43939024195762
$ws.Range("H13").Value = 9.320649069117092
$ws.Range("I13").Value = 12.10629914842097
$ws.Range("M13").Value = 73.30565342710575
$ws.Range("O13").Value = 11.93605728872384
$ws.Range("C14").Value = 4.521778523543456
$ws.Range("D14").Value = 4.416034453183623
$ws.Range("E14").Value = 51.71869852158455
$ws.Range("F14").Value = 14.08029448310079
$ws.Range("G14").Value = 12.41887220335339
$ws.Range("H14").Value = 9.343758192904097
$ws.Range("I14").Value = 12.12854290741399
$ws.Range("M14").Value = 72.69672348813351
$ws.Range("O14").Value = 11.9698662218742
$ws.Range("C15").Value = 4.505704344592509
$ws.Range("D15").Value = 4.402636420840919
$ws.Range("E15").Value = 51.43408303023234
$ws.Range("F15").Value = 14.08703510354367
$ws.Range("G15").Value = 12.40711830762945
$ws.Range("H15").Value = 9.358006082369526
$ws.Range("I15").Value = 12.14241136514456
$ws.Range("M15").Value = 72.32128940852843
$ws.Range("O15").Value = 11.99083684991599
$ws.Range("C16").Value = 4.412396617106815
$ws.Range("D16").Value = 4.32495559999153
$ws.Range("E16").Value = 49.77467610133067
$ws.Range("F16").Value = 14.13133211169708
$ws.Range("G16").Value = 12.3520980445059
$ws.Range("H16").Value = 9.441066546272076
$ws.Range("I16").Value = 12.22555299869022
$ws.Range("M16").Value = 70.13220277448985
$ws.Range("O16").Value = 12.11494097180642
$ws.Range("C17").Value = 4.354109525255939
$ws.Range("D17").Value = 4.276508791304942
$ws.Range("E17").Value = 48.73120218939457
$ws.Range("F17").Value = 14.16346129986302
$ws.Range("G17").Value = 12.32918218982192
$ws.Range("H17").Value = 9.493268293644826
$ws.Range("I17").Value = 12.27975988890409
$ws.Range("M17").Value = 68.75556782396845
$ws.Range("O17").Value = 12.1944958249117
$ws.Range("C18").Value = 4.320207874811309
$ws.Range("D18").Value = 4.248357729307262
$ws.Range("E18").Value = 48.12160123404342
$ws.Range("F18").Value = 14.18372019023938
$ws.Range("G18").Value = 12.31990890362691
$ws.Range("H18").Value = 9.523746689031791
$ws.Range("I18").Value = 12.31208768271913
$ws.Range("M18").Value = 67.95133038828438
$ws.Range("O18").Value = 12.24147778449519
$ws.Range("C19").Value = 4.308665489950588
$ws.Range("D19").Value = 4.238777723879632
$ws.Range("E19").Value = 47.91357709163687
$ws.Range("F19").Value = 14.19088107194834
$ws.Range("G19").Value = 12.31743237472913
$ws.Range("H19").Value = 9.534143562052993
$ws.Range("I19").Value = 12.32322817431193
$ws.Range("M19").Value = 67.67688957482844
$ws.Range("O19").Value = 12.25759225557313
$ws.Range("C20").Value = 4.360353432289657
$ws.Range("D20").Value = 4.281695735475675
$ws.Range("E20").Value = 48.84325565085075
$ws.Range("F20").Value = 14.15985602511353
$ws.Range("G20").Value = 12.33121540951811
$ws.Range("H20").Value = 9.487664303090478
$ws.Range("I20").Value = 12.27386990125632
$ws.Range("M20").Value = 68.90339872085933
$ws.Range("O20").Value = 12.18589965248264
$ws.Range("C21").Value = 4.529465171961054
$ws.Range("D21").Value = 4.422443044012199
$ws.Range("E21").Value = 51.85467927975165
$ws.Range("F21").Value = 14.07716391364231
$ws.Range("G21").Value = 12.42472890355965
$ws.Range("H21").Value = 9.336951096657323
$ws.Range("I21").Value = 12.12195847078695
$ws.Range("M21").Value = 72.876090609252
$ws.Range("O21").Value = 11.95988098805931
$ws.Range("C22").Value = 4.636478012592539
$ws.Range("D22").Value = 4.511776532061114
$ws.Range("E22").Value = 53.7401595660595
$ws.Range("F22").Value = 14.0398407229994
$ws.Range("G22").Value = 12.52208696136718
$ws.Range("H22").Value = 9.242595193615182
$ws.Range("I22").Value = 12.03350032155431
$ws.Range("M22").Value = 75.36282282232277
$ws.Range("O22").Value = 11.8237864535276
$ws.Range("C23").Value = 4.579691586882172
$ws.Range("D23").Value = 4.464345125165941
$ws.Range("E23").Value = 52.7413451088992
$ws.Range("F23").Value = 14.05818894464133
$ws.Range("G23").Value = 12.46675179258076
$ws.Range("H23").Value = 9.292570082991027
$ws.Range("I23").Value = 12.07969247590287
$ws.Range("M23").Value = 74.04558537521504
$ws.Range("O23").Value = 11.89532286153417
$ws.Range("C24").Value = 4.357531779888063
$ws.Range("D24").Value = 4.279351646033274
$ws.Range("E24").Value = 48.79262648894024
$ws.Range("F24").Value = 14.16148041233488
$ws.Range("G24").Value = 12.33028405686239
$ws.Range("H24").Value = 9.490196414032617
$ws.Range("I24").Value = 12.27652914698713
$ws.Range("M24").Value = 68.83660422696192
$ws.Range("O24").Value = 12.18978211251811
$ws.Range("C25").Value = 4.103633258792668
$ws.Range("D25").Value = 4.068969188123933
$ws.Range("E25").Value = 44.17236583867184
$ws.Range("F25").Value = 14.34036964180209
$ws.Range("G25").Value = 12.33102470443367
$ws.Range("H25").Value = 9.720649617242154
$ws.Range("I25").Value = 12.53247634417862
$ws.Range("M25").Value = 62.74205057483686
$ws.Range("O25").Value = 12.55384705179406
